$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the dummy variable columns
$ws.Range("D1").Value = "promocionado_25"
$ws.Range("E1").Value = "Black_Friday"

# Updated data: dates shifted forward by 4 days, new Preco values,
# and new dummy columns (promocionado_25, Black_Friday) populated.
$data = @(
    @(45918, 7172, 379.9, 1, 0),
    @(45919, 7172, 375,   1, 0),
    @(45920, 7172, 370.9, 1, 0),
    @(45921, 7172, 372,   1, 0),
    @(45922, 7172, 350,   1, 0),
    @(45923, 7172, 350,   1, 0),
    @(45924, 7172, 300,   1, 0),
    @(45925, 7172, 380,   1, 0)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Clear the lingering selection anchored on the old last cell (C9)
$ws.Range("A1").Select()

# Print setup matching the committed sheet (Letter/A4 "9" = A4, portrait)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
